$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new day column V (05-jul) with header + values
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Range("V1").Value = "05-jul"
$wsPrix.Range("U1").Copy() | Out-Null
$wsPrix.Range("V1").PasteSpecial(-4122) | Out-Null

$prixValues = @(41.53, 36.29, 37.33, 30.42, 30.27, 28.73, 28.08, 39.34, 17.86, 3.34, 0, -0.01, -0.02, -0.1, -0.11, -0.02, 0.37, 6.2, 41.25, 75.04000000000001, 67.67, 72.73999999999999, 100.32, 94.02)

for ($i = 0; $i -lt $prixValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 22).Value = $prixValues[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append a new row 19 (2025-07-03)
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
# Force text so the date-like string isn't auto-converted to a date serial,
# then clear the temporary number format so no extra style sticks around.
$wsGaz.Range("A19").NumberFormat = "@"
$wsGaz.Range("A19").Value = "2025-07-03"
$wsGaz.Range("A19").ClearFormats()
$wsGaz.Range("B19").Value = 32.85

# ---------------------------------------------------------------------------
# Sheet "CO2": append a new row 19 (2025-07-03)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A19").NumberFormat = "@"
$wsCo2.Range("A19").Value = "2025-07-03"
$wsCo2.Range("A19").ClearFormats()
$wsCo2.Range("B19").Value = 71.81
